$wb = $excel.ActiveWorkbook

# "survey" sheet: remove the plot_id row (type=integer, name=plot_id)
$survey = $wb.Worksheets.Item(1)
$survey.Rows.Item(2).Delete()

# "settings" sheet: add a new setting row table_id -> plot
$settings = $wb.Worksheets.Item(3)
$settings.Cells.Item(6, 1).Value = "table_id"
$settings.Cells.Item(6, 2).Value = "plot"

# Restore/update the selected cell on each sheet
[void]$survey.Range("B7").Select()

$choices = $wb.Worksheets.Item(2)
[void]$choices.Range("C6").Select()

[void]$settings.Range("B7").Select()
